# Remove the "Appendix: Quick prototype" section in its entirety:
#   - the "Appendix: Quick prototype" heading paragraph
#   - the blank paragraph that follows it
#   - the "Figure: PDF page 1" paragraph
#   - the paragraph containing the embedded prototype screenshot
#
# The surrounding paragraphs (the blank separator paragraph before the
# heading, and the following "Appendix: Links" heading) are left intact.

$d = $word.ActiveDocument

# Find the start of the block: the "Appendix: Quick prototype" heading.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Appendix: Quick prototype*") {
        $startPara = $p
        break
    }
}

# Find the end boundary: the next "Appendix: Links" Heading 2 paragraph
# that comes after the start paragraph (deletion stops right before it).
$endPara = $null
$passedStart = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $startPara.Range.Start) {
        $passedStart = $true
    }
    if ($passedStart -and $p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text -like "Appendix: Links*") {
        $endPara = $p
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.Start)
    $r.Delete()
}
